$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the paragraph that holds "Version 11.07.05, 2015-12-07"
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Version*") {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start
$fullText = $target.Range.Text

# The release date's day-of-month ("07") sits right after "12-" in the
# "...2015-12-07" tail; compute its offset so we don't disturb the
# "11.07.05" version number earlier in the same line.
$dayOffset = $fullText.IndexOf("12-07") + 3

$digit1Start = $pStart + $dayOffset        # position of the '0'
$digit2Start = $digit1Start + 1            # position of the '7'
$digit2End   = $digit2Start + 1

# ------------------------------------------------------------------
# 2) Drop the old _GoBack bookmark (it currently wraps the OWF logo).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3) Change the release day from 07 to 09, keeping the leading "0" in
#    its existing run and letting the new "9" land in a run of its
#    own (a temporary bookmark at the split point pins the boundary
#    so the two characters are not silently re-merged).
# ------------------------------------------------------------------
$splitRange = $d.Range($digit2Start, $digit2Start)
$d.Bookmarks.Add("ZZZ_split_marker", $splitRange)

$digitRange = $d.Range($digit2Start, $digit2End)
$digitRange.Text = "9"

$d.Bookmarks.Item("ZZZ_split_marker").Delete()

# ------------------------------------------------------------------
# 4) Re-create _GoBack right after the updated date, before the page
#    break that ends the paragraph.
# ------------------------------------------------------------------
$newBookmarkPos = $digit2Start + 1
$newBookmarkRange = $d.Range($newBookmarkPos, $newBookmarkPos)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)

Write-Output $target.Range.Text
